$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1175510589472074
$ws.Range("C2").Value = 0.9978526205381139
$ws.Range("D2").Value = 0.2444580670916839
$ws.Range("F2").Value = "Pipeline(steps=[('model', AdaBoostRegressor(n_estimators=100))])"
$ws.Range("G2").Value = 0.129829331083359
$ws.Range("H2").Value = 0.99
